# Generate Report for Handback
# The file "e1f29d49-765f-4abe-8a5d-7a268dab63cd.md" has now been
# successfully handed back (in sync with en-US) for both zh-cn and de-de.
# Update the Overview sheet plus the per-locale detail sheets accordingly,
# and clear the stale "handback not latest" error + refresh the handback
# timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the e1f29d49-... file ---
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
# row 2 = a7560759-... file, row 3 = e1f29d49-... file
$wsZhCn.Range("K2").Value = "2016-09-06 07:50:04"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-06 07:50:04"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet ---
$wsDeDe.Range("K2").Value = "2016-09-06 07:50:33"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-06 07:50:33"
$wsDeDe.Range("P3").Value = ""

# Narrow the now-empty "Error Detail" column back down on both detail sheets.
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470531463623
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470531463623
